# --------------------------------------------------------------------------
# "first of many passes through in cleaning database"
#
#  1. Shared-string text fixes:  RiboPure0.5X -> RiboPure0.5x
#                                 RiboPure0.25X -> RiboPure0.25x
#  2. Columns H / I / K (rows 2-37) get explicit TRUE()/FALSE() formulas
#     instead of bare boolean literals (values themselves are unchanged).
#  3. Sheet view: freeze the header row, scroll so the selection sits on
#     F1 / G41.
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix the RiboPure case typo everywhere it shows up in column G
# ---------------------------------------------------------------------
for ($row = 14; $row -le 25; $row++) {
    $ws.Range("G$row").Value = "RiboPure0.5x"
}
for ($row = 26; $row -le 37; $row++) {
    $ws.Range("G$row").Value = "RiboPure0.25x"
}

# ---------------------------------------------------------------------
# 2. Re-write the boolean cells in H / I / K as real formulas
#    H is FALSE() on every data row; I and K are TRUE() except for the
#    rows noted below, which keep their original FALSE value.
# ---------------------------------------------------------------------
$iFalseRows = @(11)
$kFalseRows = @(23, 30, 31, 32, 33)

for ($row = 2; $row -le 37; $row++) {
    $ws.Range("H$row").Formula = "=FALSE()"

    if ($iFalseRows -contains $row) {
        $ws.Range("I$row").Formula = "=FALSE()"
    } else {
        $ws.Range("I$row").Formula = "=TRUE()"
    }

    if ($kFalseRows -contains $row) {
        $ws.Range("K$row").Formula = "=FALSE()"
    } else {
        $ws.Range("K$row").Formula = "=TRUE()"
    }
}

# ---------------------------------------------------------------------
# 3. Sheet view: freeze row 1, and leave the view scrolled/selected the
#    way the saved workbook shows it (topLeft pane on F1, bottom-left
#    pane active cell G41).
# ---------------------------------------------------------------------
$window = $excel.ActiveWindow

$ws.Range("A2").Select()
$window.FreezePanes = $true

$ws.Range("F1").Select()
$ws.Range("G41").Select()
